$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 112 ---
# Column A (date, style matches the other date cells in the column, e.g. A111)
$ws.Range("A111").Copy($ws.Range("A112"))
$ws.Range("A112").Value = 45503.2916666667

$ws.Range("B112").Value = 0
$ws.Range("C112").Value = 3.11999988555908
$ws.Range("D112").Value = 3.11999988555908
$ws.Range("E112").Value = 3.11999988555908
$ws.Range("F112").Value = 3.11999988555908

# Column G holds the close price repeated as TEXT (shared string), not a number.
# Stage it through a formula cell so the literal keeps its textual type, then
# paste only the resulting value (not the formula) into the target cell.
$ws.Range("Z1").Formula = '="3.11999988555908"'
$ws.Range("Z1").Copy()
$ws.Range("G112").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").ClearContents()

$ws.Range("H112").Value = "ESPE.MI"

# --- Row 113 ---
$ws.Range("A111").Copy($ws.Range("A113"))
$ws.Range("A113").Value = 45504.6494444444

$ws.Range("B113").Value = 7500
$ws.Range("C113").Value = 3.10999989509583
$ws.Range("D113").Value = 3
$ws.Range("E113").Value = 3.05999994277954
$ws.Range("F113").Value = 3

$ws.Range("Z1").Formula = '="3"'
$ws.Range("Z1").Copy()
$ws.Range("G113").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").ClearContents()

$ws.Range("H113").Value = "ESPE.MI"
